# Swap the data contained in rows 2 and 3 (columns A through I) on the
# active sheet. Using Range.Copy (rather than reading/writing .Value /
# .Value2) ensures that each cell's original data type is preserved -
# in particular the "ion_id" column (A) holds numeric-looking values
# that are actually stored as text, and a naive Value2 round-trip would
# silently re-type them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = $ws.Range("A2:I2")
$row3 = $ws.Range("A3:I3")
# Scratch area well outside the used range (sheet data only goes to I25).
$scratch = $ws.Range("A100:I100")

$row2.Copy($scratch)
$row3.Copy($row2)
$scratch.Copy($row3)
$scratch.ClearContents()
